$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) "09/25/2016 (2.5hr)" -> "09/25/2016 (5" + "hr)" (two separate bold runs)
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("09/25/2016 (2.5hr)", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
$r.Text = "09/25/2016 (5"
$splitPos = $r.End
$afterRange = $d.Range($splitPos, $splitPos)
$afterRange.InsertAfter("hr)")
# Force the newly typed text into its own run (rather than being coalesced
# with the preceding run) by toggling Bold through a real transition.
$newRunRange = $d.Range($splitPos, $splitPos + 3)
$newRunRange.Bold = 0
$newRunRange.Bold = 1

# ---------------------------------------------------------------------------
# 2) Remove the trailing double space after "...during restoration." and
#    append two new bullet paragraphs, the second one holding the
#    relocated "_GoBack" bookmark.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$issue2Para = $d.Paragraphs.Item($d.Paragraphs.Count)
$trimRange = $d.Range($issue2Para.Range.End - 3, $issue2Para.Range.End - 1)
$trimRange.Text = ""

$issue2Para = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPos = $issue2Para.Range.End
$insertRange = $d.Range($insertPos, $insertPos)

$newParasXml = @"
<w:p $wNs><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Combined serialization with the tournament overall and setup, have a full human side of the game.</w:t></w:r></w:p><w:p $wNs><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t xml:space="preserve">Issue: The </w:t></w:r><w:r><w:t>reading from</w:t></w:r><w:r><w:t xml:space="preserve"> file itself seems to be fine, however, there seems to be some issue while setting the flags. As a result, while printing the game board from the tournament, some of the dice show incorrect values for top-right.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>
"@

$insertRange.InsertXML($newParasXml)
